$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    $sheetName = $ws.Name

    # ------------------------------------------------------------------
    # 1) Drop the now-unused "Mar/24" column (G) entirely. This shifts
    #    the helper table (old J/K -> I/J) left by one column and shrinks
    #    the two header merges (A2:G2 -> A2:F2, A4:G4 -> A4:F4)
    #    automatically, without touching any cell styles.
    # ------------------------------------------------------------------
    $ws.Columns.Item(7).Delete()

    # Column A is slightly wider in the new layout.
    $ws.Columns.Item(1).ColumnWidth = 14.1

    # ------------------------------------------------------------------
    # 2) Row 1 header: last month column now reads "H.mês" (was "Mar/24").
    # ------------------------------------------------------------------
    $ws.Range("F1").Value = "H.mês"

    # ------------------------------------------------------------------
    # 3) Row 2: category renamed, helper-table header moves to I2:J2.
    # ------------------------------------------------------------------
    $ws.Range("A2").Value = "Coordenação"
    $ws.Range("I2").Value = "Mês"
    $ws.Range("J2").Value = "Total Decimal"

    # ------------------------------------------------------------------
    # 4) Row 3: first staff line becomes "Eng. Júnior" with new hours.
    # ------------------------------------------------------------------
    $ws.Range("A3").Value = "Coordenação"
    $ws.Range("B3").Value = "Eng. Júnior"
    $ws.Range("D3").Value = "0,38"
    $ws.Range("E3").Value = "0,38"
    $ws.Range("F3").Value = "0,75"
    $ws.Range("I3").Value = "Jan/24"
    $ws.Range("J3").Value = 1.52

    # ------------------------------------------------------------------
    # 5) Row 4 used to be a merged "Estruturas" section header; it is now
    #    a normal staff row ("Eng. Pleno"). Unmerge, copy the row-3
    #    formatting onto it (keeps the same style index, s="3"), then
    #    fill in the values.
    # ------------------------------------------------------------------
    $ws.Range("A4:F4").UnMerge()
    $ws.Range("A3:F3").Copy()
    $ws.Range("A4:F4").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Range("A4").Value = "Coordenação"
    $ws.Range("B4").Value = "Eng. Pleno"
    $ws.Range("C4").Value = "teste1"
    $ws.Range("D4").Value = "0,38"
    $ws.Range("E4").Value = "0,38"
    $ws.Range("F4").Value = "0,75"
    $ws.Range("I4").Value = "Feb/24"
    $ws.Range("J4").Value = 1.52

    # ------------------------------------------------------------------
    # 6) Row 5: now "Eng. Sênior"; the helper-table no longer extends
    #    this far, so I5/J5 stay empty.
    # ------------------------------------------------------------------
    $ws.Range("A5").Value = "Coordenação"
    $ws.Range("B5").Value = "Eng. Sênior"
    $ws.Range("C5").Value = "teste2"
    $ws.Range("D5").Value = "0,38"
    $ws.Range("E5").Value = "0,38"
    $ws.Range("F5").Value = "0,75"
    $ws.Range("I5").ClearContents()
    $ws.Range("J5").ClearContents()

    # ------------------------------------------------------------------
    # 7) Row 6 is brand new: "Estagiário/Projetista" staff line, copying
    #    row 5's formatting.
    # ------------------------------------------------------------------
    $ws.Range("A5:F5").Copy()
    $ws.Range("A6:F6").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Range("A6").Value = "Coordenação"
    $ws.Range("B6").Value = "Estagiário/Projetista"
    $ws.Range("C6").Value = "teste3"
    $ws.Range("D6").Value = "0,38"
    $ws.Range("E6").Value = "0,38"
    $ws.Range("F6").Value = "0,75"

    # ------------------------------------------------------------------
    # 8) Move the chart down by one row and repoint its series at the
    #    relocated helper table.
    # ------------------------------------------------------------------
    $chartObj = $ws.ChartObjects(1)
    $chartObj.Top = $ws.Rows.Item(11).Top

    $chart = $chartObj.Chart
    $chart.SeriesCollection(1).Formula = "=SERIES('$sheetName'!J2,'$sheetName'!`$I`$3:`$I`$4,'$sheetName'!`$J`$3:`$J`$4,1)"
}
